# Evaluation workbook update:
#  - Add histogram ("_hist") results column (D) next to existing "_imadj" results on Sheet1
#  - Add two new rows (u_noEdge_imadj / new header row) to the existing table
#  - Add new rows 16-26 for the "_hist" variants with their own notes
#  - Add a new Sheet2 listing just the "_hist" image names
#  - Resize columns B/C and refresh the selection/active cell

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet
$ws1.Name = "Sheet1"

# ---------------------------------------------------------------------------
# Sheet1: column widths for the newly-used columns B and C
# ---------------------------------------------------------------------------
$ws1.Columns.Item(2).ColumnWidth = 27.833333333333332
$ws1.Columns.Item(3).ColumnWidth = 21.166666666666668

# ---------------------------------------------------------------------------
# Sheet1: new header row 2 (was previously part of row 3) + new row3 sub headers
# ---------------------------------------------------------------------------
$ws1.Range("B2").Value = "Processes that were effective in finding the air bubble"
$ws1.Range("B2").Font.Bold = $true

$ws1.Range("B3").Value = "no_edge_imadj"
$ws1.Range("D3").Value = "_noEdge_hist"

# ---------------------------------------------------------------------------
# Sheet1: fill in column D ("_hist" notes) for the existing imadj rows (4-13)
# ---------------------------------------------------------------------------
$ws1.Range("D4").Value  = "nothing"
$ws1.Range("D5").Value  = "nothing really"
$ws1.Range("D6").Value  = "Same as 12"
$ws1.Range("D7").Value  = "Threshold OC by recon kind of"
$ws1.Range("D8").Value  = "Regional maxima kind of"
$ws1.Range("D9").Value  = "Threshold OC by recon kind of"
$ws1.Range("D10").Value = "nothing"
$ws1.Range("D11").Value = "nothing"
$ws1.Range("D12").Value = "Watershed kinda"
$ws1.Range("D13").Value = "Threshold OC by recon"

# ---------------------------------------------------------------------------
# Sheet1: new row 14 - u_noEdge_imadj
# ---------------------------------------------------------------------------
$ws1.Range("A14").Value = "u_noEdge_imadj"
$ws1.Range("B14").Value = "Threshold OC by recon kind of"
$ws1.Range("D14").Value = "nothing"

# ---------------------------------------------------------------------------
# Sheet1: new block, rows 16-26 - the "_hist" results table
# ---------------------------------------------------------------------------
$ws1.Range("A16").Value = "2_noEdge_hist"
$ws1.Range("B16").Value = "nothing"

$ws1.Range("A17").Value = "3_noEdge_hist"
$ws1.Range("B17").Value = "nothing really"

$ws1.Range("A18").Value = "11_noEdge_hist"
$ws1.Range("B18").Value = "Same as 12"

$ws1.Range("A19").Value = "12_noEdge_hist"
$ws1.Range("B19").Value = "Threshold OC by recon kind of"

$ws1.Range("A20").Value = "e_noEdge_hist"
$ws1.Range("B20").Value = "Regional maxima kind of"
$ws1.Range("F20").Value = "can remove white spots??"

$ws1.Range("A21").Value = "f_noEdge_hist"
$ws1.Range("B21").Value = "Threshold OC by recon kind of"
$ws1.Range("C21").Value = "Regional maxima kind of"

$ws1.Range("A22").Value = "g_noEdge_hist"
$ws1.Range("B22").Value = "nothing"

$ws1.Range("A23").Value = "gg_noEdge_hist"
$ws1.Range("B23").Value = "nothing"

$ws1.Range("A24").Value = "k_noEdge_hist"
$ws1.Range("B24").Value = "Watershed kinda"
$ws1.Range("C24").Value = "Threshold OC recon "

$ws1.Range("A25").Value = "t_noEdge_hist"
$ws1.Range("B25").Value = "Threshold OC by recon"
$ws1.Range("C25").Value = [char]8730

$ws1.Range("A26").Value = "u_noEdge_hist"
$ws1.Range("B26").Value = "nothing"

# ---------------------------------------------------------------------------
# Sheet2: plain list of the "_hist" image names
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value  = "2_noEdge_hist"
$ws2.Range("A2").Value  = "3_noEdge_hist"
$ws2.Range("A3").Value  = "11_noEdge_hist"
$ws2.Range("A4").Value  = "12_noEdge_hist"
$ws2.Range("A5").Value  = "e_noEdge_hist"
$ws2.Range("A6").Value  = "f_noEdge_hist"
$ws2.Range("A7").Value  = "g_noEdge_hist"
$ws2.Range("A8").Value  = "gg_noEdge_hist"
$ws2.Range("A9").Value  = "k_noEdge_hist"
$ws2.Range("A10").Value = "t_noEdge_hist"
$ws2.Range("A11").Value = "u_noEdge_hist"

$null = $ws2.Range("A1:A11").Select()

# ---------------------------------------------------------------------------
# Re-select Sheet1 and restore the active cell / tab selection
# ---------------------------------------------------------------------------
$null = $ws1.Select()
$null = $ws1.Range("C19").Select()
